# Applies the "Updated LCA process names" edit described by the commit diff:
#  - Renames the worksheet tab from "Sheet6" to "Sheet2"
#  - Renames the 14 "Turning ..." process-contribution columns/rows to the
#    generic "Turning 1".."Turning 14" naming scheme (plus Drilling/Milling/
#    Surface Grinding/Induction Hardening), matches each process row to its
#    new column position, updates the "Total"/per-process amounts and units,
#    flips "Default units: " to "Yes", and refreshes the report Date/Time stamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet tab (Sheet6 -> Sheet2) ---
$ws.Name = "Sheet2"

# --- Report header (run Date/Time stamp) ---
$ws.Range("D1").Value2 = 45574
$ws.Range("F1").Value2 = 0.919046898148148

# --- "Default units: " flag (No -> Yes) ---
$ws.Range("B10").Value2 = "Yes"

# --- Header row 16: process/column names ---
$ws.Range("G16").Value2 = "Turning 1"
$ws.Range("H16").Value2 = "Turning 2"
$ws.Range("I16").Value2 = "Turning 3"
$ws.Range("J16").Value2 = "Turning 4"
$ws.Range("K16").Value2 = "Turning 5"
$ws.Range("L16").Value2 = "Turning 6"
$ws.Range("M16").Value2 = "Turning 7"
$ws.Range("N16").Value2 = "Turning 8"
$ws.Range("O16").Value2 = "Turning 9"
$ws.Range("P16").Value2 = "Turning 10"
$ws.Range("Q16").Value2 = "Turning 11"
$ws.Range("R16").Value2 = "Turning 12"
$ws.Range("S16").Value2 = "Turning 13"
$ws.Range("T16").Value2 = "Drilling"
$ws.Range("U16").Value2 = "Milling"
$ws.Range("V16").Value2 = "Turning 14"
$ws.Range("W16").Value2 = "Surface Grinding"
$ws.Range("X16").Value2 = "Induction Hardening"

# --- Row 17: process name, unit and amount columns ---
$ws.Range("B17").Value2 = "Drilling"
$ws.Range("D17").Value2 = "kg"
$ws.Range("E17").Value2 = 0.09132
$ws.Range("T17").Value2 = 0.09132

# --- Row 18: process name, unit and amount columns ---
$ws.Range("B18").Value2 = "Induction Hardening"
$ws.Range("D18").Value2 = "MJ"
$ws.Range("E18").Value2 = 252
$ws.Range("W18").Value2 = 0
$ws.Range("X18").Value2 = 252

# --- Row 19: process name, unit and amount columns ---
$ws.Range("B19").Value2 = "Milling"
$ws.Range("D19").Value2 = "kg"
$ws.Range("E19").Value2 = 0.02224
$ws.Range("U19").Value2 = 0.02224
$ws.Range("X19").Value2 = 0

# --- Row 20: process name, unit and amount columns ---
$ws.Range("B20").Value2 = "Primary Production Steel Billet"
$ws.Range("D20").Value2 = "kg"
$ws.Range("E20").Value2 = 13.59001
$ws.Range("F20").Value2 = 13.59001
$ws.Range("U20").Value2 = 0

# --- Row 21: process name, unit and amount columns ---
$ws.Range("B21").Value2 = "Surface Grinding"
$ws.Range("E21").Value2 = 0.001
$ws.Range("F21").Value2 = 0
$ws.Range("W21").Value2 = 0.001

# --- Row 22: process name, unit and amount columns ---
$ws.Range("B22").Value2 = "Turning 1"
$ws.Range("D22").Value2 = "kg"
$ws.Range("E22").Value2 = 0.70598
$ws.Range("G22").Value2 = 0.70598
$ws.Range("V22").Value2 = 0

# --- Row 23: process name, unit and amount columns ---
$ws.Range("B23").Value2 = "Turning 10"
$ws.Range("D23").Value2 = "kg"
$ws.Range("E23").Value2 = 0.07568
$ws.Range("P23").Value2 = 0.07568
$ws.Range("R23").Value2 = 0

# --- Row 24: process name, unit and amount columns ---
$ws.Range("B24").Value2 = "Turning 11"
$ws.Range("D24").Value2 = "kg"
$ws.Range("E24").Value2 = 0.78442
$ws.Range("G24").Value2 = 0
$ws.Range("Q24").Value2 = 0.78442

# --- Row 25: process name, unit and amount columns ---
$ws.Range("B25").Value2 = "Turning 12"
$ws.Range("E25").Value2 = 0.00264
$ws.Range("K25").Value2 = 0
$ws.Range("R25").Value2 = 0.00264

# --- Row 26: process name, unit and amount columns ---
$ws.Range("B26").Value2 = "Turning 13"
$ws.Range("D26").Value2 = "kg"
$ws.Range("E26").Value2 = 0.00445
$ws.Range("S26").Value2 = 0.00445

# --- Row 27: process name, unit and amount columns ---
$ws.Range("B27").Value2 = "Turning 14"
$ws.Range("D27").Value2 = "kg"
$ws.Range("E27").Value2 = 0.00043
$ws.Range("J27").Value2 = 0
$ws.Range("V27").Value2 = 0.00043

# --- Row 28: process name, unit and amount columns ---
$ws.Range("B28").Value2 = "Turning 2"
$ws.Range("D28").Value2 = "kg"
$ws.Range("E28").Value2 = 0.21985
$ws.Range("H28").Value2 = 0.21985
$ws.Range("M28").Value2 = 0

# --- Row 29: process name, unit and amount columns ---
$ws.Range("B29").Value2 = "Turning 3"
$ws.Range("D29").Value2 = "kg"
$ws.Range("E29").Value2 = 0.00043
$ws.Range("I29").Value2 = 0.00043

# --- Row 30: process name, unit and amount columns ---
$ws.Range("B30").Value2 = "Turning 4"
$ws.Range("D30").Value2 = "kg"
$ws.Range("E30").Value2 = 0.04842
$ws.Range("H30").Value2 = 0
$ws.Range("J30").Value2 = 0.04842

# --- Row 31: process name, unit and amount columns ---
$ws.Range("B31").Value2 = "Turning 5"
$ws.Range("E31").Value2 = 2.45273
$ws.Range("K31").Value2 = 2.45273
$ws.Range("L31").Value2 = 0

# --- Row 32: process name, unit and amount columns ---
$ws.Range("B32").Value2 = "Turning 6"
$ws.Range("D32").Value2 = "kg"
$ws.Range("E32").Value2 = 4.8298
$ws.Range("L32").Value2 = 4.8298
$ws.Range("N32").Value2 = 0

# --- Row 33: process name, unit and amount columns ---
$ws.Range("B33").Value2 = "Turning 7"
$ws.Range("D33").Value2 = "kg"
$ws.Range("E33").Value2 = 0.00316
$ws.Range("M33").Value2 = 0.00316
$ws.Range("O33").Value2 = 0

# --- Row 34: process name, unit and amount columns ---
$ws.Range("B34").Value2 = "Turning 8"
$ws.Range("D34").Value2 = "kg"
$ws.Range("E34").Value2 = 0.75184
$ws.Range("N34").Value2 = 0.75184
$ws.Range("P34").Value2 = 0

# --- Row 35: process name, unit and amount columns ---
$ws.Range("B35").Value2 = "Turning 9"
$ws.Range("D35").Value2 = "kg"
$ws.Range("E35").Value2 = 0.02251
$ws.Range("O35").Value2 = 0.02251
$ws.Range("Q35").Value2 = 0
